$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 8
    3  = 5
    4  = 5
    5  = 2
    6  = 5
    7  = 4
    8  = 4
    9  = 5
    10 = 6
    11 = 4
    12 = 5
    13 = 2
    14 = 5
    15 = 3
    16 = 4
    17 = 2
    18 = 4
    19 = 1
    20 = 7
    21 = 6
    22 = 8
    23 = 6
    24 = 6
    25 = 5
    26 = 7
    27 = 5
    28 = 3
    29 = 7
    30 = 4
    31 = 2
    32 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
